$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.404879212379456
$ws.Range("B1").Value = 3.935016632080078
$ws.Range("C1").Value = 5.655776023864746
$ws.Range("D1").Value = 1.614978551864624
$ws.Range("E1").Value = 0.9596564173698425
